# "make BinaryLinear easier to use"
# - Rename Sheet1 -> prob01
# - Add a new sheet "prob05" after prob01 (becomes the active tab)
# - Populate prob05 with u/v inputs and EXP/POWER based formulas

$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "prob01"

# Insert the new sheet right after prob01; it becomes the active sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "prob05"

# Headers.
$ws2.Range("A1").Value = "u"
$ws2.Range("B1").Value = "v"
$ws2.Range("C1").Value = "exp_v"
$ws2.Range("D1").Value = "exp_u"
$ws2.Range("E1").Value = "u_exp_v"
$ws2.Range("F1").Value = "2_v_exp_minus_u"
$ws2.Range("G1").Value = "E_u_v"

# Inputs.
$ws2.Range("A2").Value = 0.045
$ws2.Range("B2").Value = 0.024

# Formulas.
$ws2.Range("C2").Formula = "=EXP(`$B2)"
$ws2.Range("D2").Formula = "=EXP(`$A2)"
$ws2.Range("E2").Formula = "=`$A2*`$C2"
$ws2.Range("F2").Formula = "=2*`$B2/`$D2"
$ws2.Range("G2").Formula = "=POWER(`$E2-`$F2, 2)"

$ws2.Range("E9").Select()
